# DatabaseToExcel.xlsx edit:
#  - rename the data sheet to "Student Data"
#  - insert two new student rows (Samay Raina, Sagar Shah) above the
#    existing "Suhana Sharma" row
#  - append two more new student rows (Ellen Degenerous, Nihar) at the
#    bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet2 is blank in this workbook both before and after the edit, so it
# is safe to use a single cell on it as scratch space: writing a value
# there with a Text number format, copying it and pasting it (values +
# formats) into the real destination is the only way this COM surface
# will store a numeric-looking string ("90", "2023-12-07", ...) as a
# literal shared string instead of silently coercing it to a number /
# date serial. Clearing the scratch cell afterwards leaves Sheet2 empty
# again.
$scratch = $wb.Worksheets.Item(2)

function Set-TextValue($targetCell, $value) {
    $scratch.Cells.Item(1, 1).NumberFormat = "@"
    $scratch.Cells.Item(1, 1).Value = $value
    $scratch.Cells.Item(1, 1).Copy()
    $targetCell.PasteSpecial()
    $scratch.Cells.Item(1, 1).Clear()
}

function Set-TextRow($row, $startCol, $values) {
    $n = $values.Length
    $srcRange = $scratch.Range($scratch.Cells.Item(1, 1), $scratch.Cells.Item(1, $n))
    $srcRange.NumberFormat = "@"
    for ($i = 0; $i -lt $n; $i++) {
        $scratch.Cells.Item(1, $i + 1).Value = $values[$i]
    }
    $srcRange.Copy()
    $dstRange = $ws.Range($ws.Cells.Item($row, $startCol), $ws.Cells.Item($row, $startCol + $n - 1))
    $dstRange.PasteSpecial()
    $srcRange.Clear()
}

# 1. Rename the first sheet.
$ws.Name = "Student Data"

# 2. Make room for the two new rows right above the "Suhana Sharma" row
#    (currently row 3); this pushes it down to row 5, matching the diff.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# 3. Plain-text columns (A, B, D, F, G, H, N) for every new row can be
#    written directly - none of these values look like a number or a
#    date, so the COM layer stores them as ordinary shared strings.
$plainRows = @{
    3 = @("Samay Raina", "2023PCP5320", "niharkajla28@gmail.com", "Male", "PG", "PS", "General")
    4 = @("Sagar Shah", "2023PCP5319", "niharkajla123@gmail.com", "Male", "PG", "VLSI", "OBC")
    6 = @("Ellen Degenerous", "2023PCP5318", "niharkajla28@gmail.com", "Female", "PG", "CSIS", "General")
    7 = @("Nihar", "2023PCP5317", "niharkajla28@gmail.com", "Male", "PG", "CSE", "General")
}
$plainCols = @(1, 2, 4, 6, 7, 8, 14)

foreach ($row in @(3, 4, 6, 7)) {
    $values = $plainRows[$row]
    for ($i = 0; $i -lt $plainCols.Length; $i++) {
        $ws.Cells.Item($row, $plainCols[$i]).Value = $values[$i]
    }
}

# 4. Numeric-looking text columns: contact number (C), date of birth (E)
#    and the marks/CGPA/backlogs/red-flags block (I:M) - these need the
#    text-forcing helper above.
Set-TextValue $ws.Cells.Item(3, 3) "8050106439"
Set-TextValue $ws.Cells.Item(4, 3) "8050106439"
Set-TextValue $ws.Cells.Item(6, 3) "8050106439"
Set-TextValue $ws.Cells.Item(7, 3) "8050106439"

Set-TextValue $ws.Cells.Item(3, 5) "2023-12-07"
Set-TextValue $ws.Cells.Item(4, 5) "2024-02-01"
Set-TextValue $ws.Cells.Item(6, 5) "2024-03-05"
Set-TextValue $ws.Cells.Item(7, 5) "1995-10-28"

Set-TextRow 3 9 @("90", "95", "8", "1", "0")
Set-TextRow 4 9 @("80", "85", "8.75", "0", "0")
Set-TextRow 6 9 @("70", "70", "9", "0", "0")
Set-TextRow 7 9 @("88", "85", "8.167", "1", "2")

# 5. Column A now holds the longer name "Ellen Degenerous"; widen it to
#    keep it a best-fit-style width similar to the rest of the table.
$ws.Columns.Item(1).ColumnWidth = 15
